# The authoritative edit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml becomes the "Integral" theme (colours
# previously at theme2.xml) and theme2.xml becomes the stock "Office Theme"
# (colours previously at theme1.xml). The <a:fontScheme>/<a:fmtScheme>
# blocks of both theme parts are already byte-identical, so the only
# effective change is to the 12-slot colour scheme.
#
# This host's PowerPoint object model only exposes/持ちます a single writable
# Theme - the one wired to the presentation's slide master
# (ppt/theme/theme2.xml) - via SlideMaster.Theme.ThemeColorScheme. We drive
# that object to replace the "Integral" palette with the stock "Office"
# palette, matching the diff's target content for theme2.xml.
#
# ThemeColorScheme.Colors(i) uses the clrScheme schema order:
#   1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#   9=accent5 10=accent6 11=hlink 12=folHlink
# .RGB takes/returns a standard VBA RGB() colour (R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
